$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DeckValues text (shared across C2:C4) to the new pipe-delimited format
$newVal = "N|5,N|6,N|7,N|8,N|9,N|5,N|6,N|7,N|8,N|9"
$ws.Range("C2").Value = $newVal
$ws.Range("C3").Value = $newVal
$ws.Range("C4").Value = $newVal

# Widen column C so the longer text fits (results in a serialized width of 39)
$ws.Columns.Item(3).ColumnWidth = 38.16666666666664

# Add a new bold-styled (otherwise empty) cell at I18
$ws.Range("I18").Font.Bold = $true

# Move/refresh the active selection to C2
$ws.Range("C2").Select() | Out-Null

# Page setup tweaks
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
